$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Add a new login row (row 4) mirroring the existing rows' pattern:
# A=firmid, B=<user name>, C=password, D=newautomation
$ws.Range("A4").Value = "firmid"
$ws.Range("B4").Value = "Mohan Kumar"
$ws.Range("C4").Value = "password"
$ws.Range("D4").Value = "newautomation"

# Move the active selection to the newly added cell, matching the
# workbook's recorded UI state after the edit.
[void]$ws.Range("A4").Select()
